$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$panes = $win.Panes
Write-Host ($panes | Get-Member | Out-String)
